# Natmi following Dr Hou advice:
# Rebuild the LR-pairs table (Dhh -> Ptch1) as the full 3x3 cross of
# Sending/Target clusters (FAPs, sCs, ECs) with updated edge statistics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("FAPs", "Dhh", "Ptch1", "FAPs", 2, 0.6666666666666666, 1.853892, 5.561676, 0.5711238486747862, 0.571123848674786, 2, 0.6666666666666666, 5.793595666666666, 17.380787, 0.3620794979747173, 0.3620794979747173, 10.740700657668, 96.66630591901199, 0.206792236409555, 0.206792236409555),
    @("FAPs", "Dhh", "Ptch1", "sCs", 2, 0.6666666666666666, 1.853892, 5.561676, 0.5711238486747862, 0.571123848674786, 3, 1, 4.518509, 13.555527, 0.2823910339009808, 0.2823910339009808, 8.376827687028, 75.391449183252, 0.1612802541127802, 0.1612802541127801),
    @("FAPs", "Dhh", "Ptch1", "ECs", 2, 0.6666666666666666, 1.853892, 5.561676, 0.5711238486747862, 0.571123848674786, 3, 1, 5.688789333333333, 17.066368, 0.3555294681243019, 0.355529468124302, 10.546401034752, 94.91760931276801, 0.203051358152451, 0.203051358152451),
    @("sCs", "Dhh", "Ptch1", "FAPs", 3, 1, 0.4262446666666667, 1.278734, 0.1313121230922664, 0.1313121230922664, 2, 0.6666666666666666, 5.793595666666666, 17.380787, 0.3620794979747173, 0.3620794979747173, 2.469489253739777, 22.225403283658, 0.04754542760724211, 0.04754542760724211),
    @("sCs", "Dhh", "Ptch1", "sCs", 3, 1, 0.4262446666666667, 1.278734, 0.1313121230922664, 0.1313121230922664, 3, 1, 4.518509, 13.555527, 0.2823910339009808, 0.2823910339009808, 1.925990362535333, 17.333913262818, 0.03708136620375797, 0.03708136620375796),
    @("sCs", "Dhh", "Ptch1", "ECs", 3, 1, 0.4262446666666667, 1.278734, 0.1313121230922664, 0.1313121230922664, 3, 1, 5.688789333333333, 17.066368, 0.3555294681243019, 0.355529468124302, 2.424816113123555, 21.823345018112, 0.04668532928126634, 0.04668532928126634),
    @("ECs", "Dhh", "Ptch1", "FAPs", 3, 1, 0.9659053333333333, 2.897716, 0.2975640282329475, 0.2975640282329475, 2, 0.6666666666666666, 5.793595666666666, 17.380787, 0.3620794979747173, 0.3620794979747173, 5.596064953610221, 50.36458458249199, 0.1077418339579202, 0.1077418339579202),
    @("ECs", "Dhh", "Ptch1", "sCs", 3, 1, 0.9659053333333333, 2.897716, 0.2975640282329475, 0.2975640282329475, 3, 1, 4.518509, 13.555527, 0.2823910339009808, 0.2823910339009808, 4.364451941814666, 39.280067476332, 0.08402941358444269, 0.08402941358444269),
    @("ECs", "Dhh", "Ptch1", "ECs", 3, 1, 0.9659053333333333, 2.897716, 0.2975640282329475, 0.2975640282329475, 3, 1, 5.688789333333333, 17.066368, 0.3555294681243019, 0.355529468124302, 5.494831957276444, 49.453487615488, 0.1057927806905846, 0.1057927806905846),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $row[$j]
    }
}